# fix(publipostage): Correct status name
# Replace the "bleu" status label with "noir" and correct the
# "pas de résultat ni de publication" status name wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    if ($bCell.Value2 -eq "bleu") {
        $bCell.Value2 = "noir"
    }

    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq "pas de résultat ni de publication") {
        $cCell.Value2 = "pas de résultat postés ni publiés"
    }
}
